# "Split the assignment screen to separate wells"
#
# On the "To fix" sheet, two rows that were still marked "?" / "Open" are
# now marked "Done" (B10 and B11), and the AutoFilter is turned into an
# active filter on the Status column (col B) that only shows "?" and
# "Open" rows - hiding everything else (Done / Rejected), including the
# two rows that just became "Done".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the status of the last two bugs to "Done".
$ws.Range("B10").Value = "Done"
$ws.Range("B11").Value = "Done"

# Re-apply the autofilter over the whole table (A1:B11) filtering column B
# (the 2nd column of the range) to only the "?" and "Open" values. This
# also hides every row whose status isn't "?" or "Open".
$ws.AutoFilterMode = $false
$ws.Range("A1:B11").AutoFilter(2, @("?", "Open"), 7)

# Move the active selection to B12, right below the table.
$ws.Range("B12").Select()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# new filter range.
$n = $wb.Names.Item(1)
$n.RefersTo = "='To fix'!`$A`$1:`$B`$11"
